$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1: copy formatting from H1 (bold/border/center) then set labels
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I = innings start (I0), J = innings finish (IF)
$data = @{
    2 = @(9, 9)
    3 = @(10, 10)
    4 = @(9, 9)
    5 = @(8, 8)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(11, 11)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(10, 10)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(9, 9)
    28 = @(9, 9)
    29 = @(8, 9)
    30 = @(9, 9)
    31 = @(9, 10)
    32 = @(8, 8)
    33 = @(9, 9)
    34 = @(8, 8)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(9, 9)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(9, 9)
    45 = @(9, 9)
    46 = @(9, 10)
    47 = @(9, 9)
    48 = @(9, 9)
    49 = @(9, 9)
    50 = @(9, 9)
    51 = @(8, 9)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(9, 9)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(8, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(7, 7)
    67 = @(7, 7)
    68 = @(4, 4)
    69 = @(6, 6)
    70 = @(3, 3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Output "done"
